# Update odds/liquidity values for Jogos_do_Dia_Betfair_Back_Lay_2026-01-12.xlsx
# Values below reflect the latest Betfair Back/Lay snapshot for 2026-01-12.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.85
$ws.Range("G2").Value = 1.97
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 3.55
$ws.Range("O2").Value = 1.31
$ws.Range("R2").Value = 1.27
$ws.Range("S2").Value = 1.91
$ws.Range("V2").Value = 1.26
$ws.Range("W2").Value = 2.02
$ws.Range("Z2").Value = 980
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 980
$ws.Range("AE2").Value = 70
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 980
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 980
$ws.Range("AK2").Value = 980
$ws.Range("AL2").Value = 980
$ws.Range("AM2").Value = 1000
$ws.Range("AO2").Value = 1000

# Row 3
$ws.Range("J3").Value = 3.1
$ws.Range("Q3").Value = 2.16
$ws.Range("S3").Value = 4.1
$ws.Range("T3").Value = 1.89

# Row 4
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 2.18
$ws.Range("Q4").Value = 1.66
$ws.Range("U4").Value = 1.55

# Row 5
$ws.Range("F5").Value = 1.74
$ws.Range("G5").Value = 2.14
$ws.Range("H5").Value = 3.9
$ws.Range("I5").Value = 6.8
$ws.Range("J5").Value = 3.6
$ws.Range("K5").Value = 980
$ws.Range("N5").Value = 1.01
$ws.Range("O5").Value = 1.2
$ws.Range("P5").Value = 1.45
$ws.Range("Q5").Value = 1.2
$ws.Range("R5").Value = 1.51
$ws.Range("S5").Value = 2.34
$ws.Range("T5").Value = 1.01
$ws.Range("U5").Value = 1.01
$ws.Range("V5").Value = 1.22
$ws.Range("W5").Value = 2.06
$ws.Range("X5").Value = 32
$ws.Range("Y5").Value = 30
$ws.Range("Z5").Value = 980
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 14.5
$ws.Range("AD5").Value = 29
$ws.Range("AF5").Value = 18
$ws.Range("AG5").Value = 15.5
$ws.Range("AH5").Value = 1000
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 29
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 44
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 980

# Row 6
$ws.Range("L6").Value = 1.57
$ws.Range("U6").Value = 1.81
$ws.Range("Y6").Value = 12.5
$ws.Range("AH6").Value = 26

# Row 7
$ws.Range("F7").Value = 1.96
$ws.Range("I7").Value = 4.1
$ws.Range("J7").Value = 3.9
$ws.Range("K7").Value = 980
$ws.Range("M7").Value = 1.03
$ws.Range("N7").Value = 2.16
$ws.Range("O7").Value = 1.18
$ws.Range("P7").Value = 2.16
$ws.Range("R7").Value = 1.22
$ws.Range("S7").Value = 1.57
$ws.Range("T7").Value = 1.01
$ws.Range("U7").Value = 1.01
$ws.Range("X7").Value = 980
$ws.Range("Y7").Value = 1000
$ws.Range("Z7").Value = 1000
$ws.Range("AA7").Value = 60
$ws.Range("AB7").Value = 1000
$ws.Range("AC7").Value = 1000
$ws.Range("AD7").Value = 1000
$ws.Range("AE7").Value = 1000
$ws.Range("AF7").Value = 1000
$ws.Range("AG7").Value = 1000
$ws.Range("AH7").Value = 1000
$ws.Range("AI7").Value = 1000
$ws.Range("AJ7").Value = 980
$ws.Range("AK7").Value = 1000
$ws.Range("AL7").Value = 1000
$ws.Range("AM7").Value = 1000
$ws.Range("AN7").Value = 1000
$ws.Range("AO7").Value = 980

# Row 8
$ws.Range("F8").Value = 2.26
$ws.Range("G8").Value = 2.7
$ws.Range("H8").Value = 2.66
$ws.Range("I8").Value = 3.3
$ws.Range("J8").Value = 3.95
$ws.Range("N8").Value = 1.79
$ws.Range("P8").Value = 1.78
$ws.Range("Q8").Value = 1.31
$ws.Range("R8").Value = 1.79
$ws.Range("T8").Value = 1.01
$ws.Range("V8").Value = 1.43
$ws.Range("W8").Value = 1.63
$ws.Range("Y8").Value = 980
$ws.Range("Z8").Value = 980
$ws.Range("AB8").Value = 980
$ws.Range("AC8").Value = 1000
$ws.Range("AD8").Value = 1000
$ws.Range("AE8").Value = 1000
$ws.Range("AF8").Value = 980
$ws.Range("AG8").Value = 1000
$ws.Range("AH8").Value = 1000
$ws.Range("AI8").Value = 980
$ws.Range("AK8").Value = 980
$ws.Range("AL8").Value = 980

# Row 9
$ws.Range("F9").Value = 2.08
$ws.Range("G9").Value = 2.6
$ws.Range("H9").Value = 2.9
$ws.Range("I9").Value = 3.75
$ws.Range("J9").Value = 3.6
$ws.Range("K9").Value = 5
$ws.Range("M9").Value = 1.02
$ws.Range("O9").Value = 1.12
$ws.Range("Q9").Value = 1.36
$ws.Range("R9").Value = 1.74
$ws.Range("S9").Value = 1.87
$ws.Range("W9").Value = 1.69
$ws.Range("X9").Value = 50
$ws.Range("Y9").Value = 34
$ws.Range("Z9").Value = 44
$ws.Range("AB9").Value = 29
$ws.Range("AD9").Value = 22
$ws.Range("AE9").Value = 40
$ws.Range("AF9").Value = 30
$ws.Range("AH9").Value = 21
$ws.Range("AI9").Value = 44
$ws.Range("AJ9").Value = 44
$ws.Range("AK9").Value = 29
$ws.Range("AL9").Value = 34
$ws.Range("AM9").Value = 65

# Row 10
$ws.Range("F10").Value = 1.79
$ws.Range("G10").Value = 1.97
$ws.Range("H10").Value = 3.75
$ws.Range("I10").Value = 4.5
$ws.Range("M10").Value = 1.02
$ws.Range("N10").Value = 6.6
$ws.Range("O10").Value = 1.14
$ws.Range("P10").Value = 2.9
$ws.Range("Q10").Value = 1.42
$ws.Range("R10").Value = 1.77
$ws.Range("S10").Value = 2.04
$ws.Range("T10").Value = 1.48
$ws.Range("U10").Value = 2.66
$ws.Range("V10").Value = 1.28
$ws.Range("W10").Value = 2.04
$ws.Range("X10").Value = 42
$ws.Range("Y10").Value = 32
$ws.Range("Z10").Value = 40
$ws.Range("AA10").Value = 80
$ws.Range("AB10").Value = 20
$ws.Range("AC10").Value = 14.5
$ws.Range("AD10").Value = 22
$ws.Range("AE10").Value = 46
$ws.Range("AF10").Value = 20
$ws.Range("AG10").Value = 14
$ws.Range("AH10").Value = 19
$ws.Range("AI10").Value = 44
$ws.Range("AJ10").Value = 27
$ws.Range("AK10").Value = 21
$ws.Range("AL10").Value = 29
$ws.Range("AM10").Value = 55
$ws.Range("AN10").Value = 8.199999999999999
$ws.Range("AO10").Value = 28

# Row 11
$ws.Range("F11").Value = 2.02
$ws.Range("G11").Value = 2.28
$ws.Range("H11").Value = 3.75
$ws.Range("L11").Value = 1.43
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 1.25
$ws.Range("P11").Value = 1.25
$ws.Range("Q11").Value = 1.36
$ws.Range("R11").Value = 1.18
$ws.Range("S11").Value = 1.36
$ws.Range("T11").Value = 1.84
$ws.Range("U11").Value = 1.01
$ws.Range("W11").Value = 1.7
$ws.Range("X11").Value = 1000
$ws.Range("Y11").Value = 1000
$ws.Range("AA11").Value = 1000
$ws.Range("AB11").Value = 1000
$ws.Range("AC11").Value = 1000
$ws.Range("AD11").Value = 1000
$ws.Range("AE11").Value = 65
$ws.Range("AF11").Value = 1000
$ws.Range("AG11").Value = 1000
$ws.Range("AI11").Value = 1000
$ws.Range("AM11").Value = 1000
$ws.Range("AN11").Value = 980
$ws.Range("AO11").Value = 1000

# Row 12
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 980
$ws.Range("H12").Value = 2.62
$ws.Range("I12").Value = 2.9
$ws.Range("J12").Value = 2.92
$ws.Range("K12").Value = 3.2
$ws.Range("L12").Value = 1.01
$ws.Range("M12").Value = 1.09
$ws.Range("N12").Value = 2.52
$ws.Range("O12").Value = 1.48
$ws.Range("P12").Value = 1.59
$ws.Range("Q12").Value = 1.48
$ws.Range("R12").Value = 1.17
$ws.Range("S12").Value = 4
$ws.Range("T12").Value = 1.8
$ws.Range("U12").Value = 1.01
$ws.Range("V12").Value = 1.52
$ws.Range("W12").Value = 1.44
$ws.Range("X12").Value = 1000
$ws.Range("Y12").Value = 1000
$ws.Range("Z12").Value = 1000
$ws.Range("AA12").Value = 1000
$ws.Range("AB12").Value = 1000
$ws.Range("AC12").Value = 1000
$ws.Range("AD12").Value = 1000
$ws.Range("AE12").Value = 1000
$ws.Range("AF12").Value = 1000
$ws.Range("AG12").Value = 1000
$ws.Range("AH12").Value = 1000
$ws.Range("AI12").Value = 1000
$ws.Range("AJ12").Value = 1000
$ws.Range("AK12").Value = 1000
$ws.Range("AL12").Value = 1000
$ws.Range("AM12").Value = 1000
$ws.Range("AN12").Value = 1000
$ws.Range("AO12").Value = 1000

# Row 13
$ws.Range("L13").Value = 1.29
$ws.Range("P13").Value = 2.36
$ws.Range("Q13").Value = 1.69
$ws.Range("U13").Value = 1.68
$ws.Range("V13").Value = 1.06
$ws.Range("W13").Value = 4.5
$ws.Range("X13").Value = 22
$ws.Range("Z13").Value = 160
$ws.Range("AA13").Value = 920
$ws.Range("AD13").Value = 55
$ws.Range("AE13").Value = 320
$ws.Range("AG13").Value = 11.5
$ws.Range("AH13").Value = 40
$ws.Range("AI13").Value = 240
$ws.Range("AK13").Value = 14.5
$ws.Range("AL13").Value = 46
$ws.Range("AM13").Value = 290
$ws.Range("AO13").Value = 470

# Row 14
$ws.Range("L14").Value = 1.44
$ws.Range("M14").Value = 1.08
$ws.Range("O14").Value = 1.38
$ws.Range("P14").Value = 1.81
$ws.Range("Q14").Value = 2.16
$ws.Range("S14").Value = 4
$ws.Range("V14").Value = 1.53
$ws.Range("W14").Value = 1.53
$ws.Range("X14").Value = 11.5
$ws.Range("Z14").Value = 17.5
$ws.Range("AA14").Value = 44
$ws.Range("AB14").Value = 11
$ws.Range("AD14").Value = 12.5
$ws.Range("AE14").Value = 32
$ws.Range("AF14").Value = 17.5
$ws.Range("AG14").Value = 12.5
$ws.Range("AI14").Value = 48
$ws.Range("AJ14").Value = 44
$ws.Range("AK14").Value = 34
$ws.Range("AL14").Value = 48
$ws.Range("AM14").Value = 110
$ws.Range("AN14").Value = 32
$ws.Range("AO14").Value = 32
